# Add data for 2024-04-07
# Updates column K (year 2024 cumulative totals) across the Citywide Totals,
# By Neighborhood, and individual neighborhood sheets to include the new week.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1861  # Aggravated Assault: 1836 -> 1861
$ws.Range("K3").Value = 1774  # Aggravated Battery: 1751 -> 1774
$ws.Range("K4").Value = 386  # Criminal Sexual Assault: 383 -> 386
$ws.Range("K5").Value = 118  # Homicide: 116 -> 118
$ws.Range("K6").Value = 2288  # Robbery: 2263 -> 2288
$ws.Range("K7").Value = 6427  # Total: 6349 -> 6427

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 21  # Aggravated Battery: 20 -> 21
$ws.Range("K6").Value = 55  # Robbery: 54 -> 55
$ws.Range("K7").Value = 101  # Total: 99 -> 101

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 123  # Aggravated Assault: 121 -> 123
$ws.Range("K3").Value = 122  # Aggravated Battery: 120 -> 122
$ws.Range("K6").Value = 144  # Robbery: 143 -> 144
$ws.Range("K7").Value = 420  # Total: 415 -> 420

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 55  # Aggravated Assault: 54 -> 55
$ws.Range("K3").Value = 42  # Aggravated Battery: 41 -> 42
$ws.Range("K7").Value = 135  # Total: 133 -> 135

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 74  # Aggravated Assault: 73 -> 74
$ws.Range("K3").Value = 98  # Aggravated Battery: 97 -> 98
$ws.Range("K7").Value = 260  # Total: 258 -> 260

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 51  # Aggravated Assault: 50 -> 51
$ws.Range("K7").Value = 212  # Total: 211 -> 212

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 39  # Aggravated Battery: 38 -> 39
$ws.Range("K7").Value = 156  # Total: 155 -> 156

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 51  # Albany Park: 50 -> 51
$ws.Range("K7").Value = 178  # Auburn Gresham: 177 -> 178
$ws.Range("K8").Value = 420  # Austin: 415 -> 420
$ws.Range("K13").Value = 10  # Boystown: 9 -> 10
$ws.Range("K14").Value = 36  # Bridgeport: 35 -> 36
$ws.Range("K17").Value = 10  # Burnside: 9 -> 10
$ws.Range("K18").Value = 47  # Calumet Heights: 46 -> 47
$ws.Range("K19").Value = 177  # Chatham: 174 -> 177
$ws.Range("K20").Value = 142  # Chicago Lawn: 139 -> 142
$ws.Range("K22").Value = 17  # Clearing: 16 -> 17
$ws.Range("K23").Value = 58  # Douglas: 57 -> 58
$ws.Range("K29").Value = 310  # Englewood: 303 -> 310
$ws.Range("K31").Value = 73  # Gage Park: 72 -> 73
$ws.Range("K33").Value = 260  # Garfield Park: 258 -> 260
$ws.Range("K35").Value = 10  # Gold Coast: 9 -> 10
$ws.Range("K37").Value = 212  # Grand Crossing: 211 -> 212
$ws.Range("K42").Value = 221  # Humboldt Park: 219 -> 221
$ws.Range("K43").Value = 62  # Hyde Park: 61 -> 62
$ws.Range("K44").Value = 62  # Irving Park: 60 -> 62
$ws.Range("K46").Value = 13  # Jefferson Park: 12 -> 13
$ws.Range("K48").Value = 77  # Lake View: 69 -> 77
$ws.Range("K49").Value = 38  # Lincoln Park: 37 -> 38
$ws.Range("K50").Value = 35  # Lincoln Square: 34 -> 35
$ws.Range("K51").Value = 75  # Little Italy, UIC: 74 -> 75
$ws.Range("K53").Value = 101  # Logan Square: 99 -> 101
$ws.Range("K54").Value = 107  # Loop: 105 -> 107
$ws.Range("K55").Value = 68  # Lower West Side: 67 -> 68
$ws.Range("K56").Value = 10  # Magnificent Mile: 9 -> 10
$ws.Range("K57").Value = 14  # Mckinley Park: 15 -> 14
$ws.Range("K60").Value = 43  # Morgan Park: 42 -> 43
$ws.Range("K64").Value = 44  # Near South Side: 42 -> 44
$ws.Range("K65").Value = 156  # New City: 155 -> 156
$ws.Range("K67").Value = 248  # North Lawndale: 247 -> 248
$ws.Range("K72").Value = 26  # Old Town: 24 -> 26
$ws.Range("K75").Value = 23  # Pullman: 22 -> 23
$ws.Range("K76").Value = 92  # River North: 91 -> 92
$ws.Range("K79").Value = 172  # Roseland: 170 -> 172
$ws.Range("K83").Value = 135  # South Chicago: 133 -> 135
$ws.Range("K85").Value = 324  # South Shore: 319 -> 324
$ws.Range("K86").Value = 44  # Streeterville: 43 -> 44
$ws.Range("K91").Value = 57  # Washington Park: 56 -> 57
$ws.Range("K94").Value = 76  # West Loop: 75 -> 76
$ws.Range("K96").Value = 92  # West Ridge: 91 -> 92
$ws.Range("K97").Value = 58  # West Town: 55 -> 58
$ws.Range("K101").Value = 6427  # Total: 6349 -> 6427

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 26  # Aggravated Assault: 25 -> 26
$ws.Range("K7").Value = 73  # Total: 72 -> 73

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 71  # Aggravated Assault: 70 -> 71
$ws.Range("K7").Value = 248  # Total: 247 -> 248

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 25  # Robbery: 24 -> 25
$ws.Range("K7").Value = 38  # Total: 37 -> 38

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 37  # Aggravated Battery: 35 -> 37
$ws.Range("K7").Value = 107  # Total: 105 -> 107

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 83  # Aggravated Assault: 81 -> 83
$ws.Range("K3").Value = 106  # Aggravated Battery: 104 -> 106
$ws.Range("K5").Value = 8  # Homicide: 7 -> 8
$ws.Range("K6").Value = 98  # Robbery: 96 -> 98
$ws.Range("K7").Value = 310  # Total: 303 -> 310

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 12  # Aggravated Battery: 11 -> 12
$ws.Range("K4").Value = 14  # Criminal Sexual Assault: 11 -> 14
$ws.Range("K6").Value = 35  # Robbery: 31 -> 35
$ws.Range("K7").Value = 77  # Total: 69 -> 77

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 55  # Aggravated Assault: 54 -> 55
$ws.Range("K3").Value = 51  # Aggravated Battery: 50 -> 51
$ws.Range("K6").Value = 58  # Robbery: 57 -> 58
$ws.Range("K7").Value = 177  # Total: 174 -> 177

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 27  # Robbery: 25 -> 27
$ws.Range("K7").Value = 62  # Total: 60 -> 62

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 18  # Aggravated Battery: 17 -> 18
$ws.Range("K7").Value = 92  # Total: 91 -> 92

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 16  # Aggravated Assault: 15 -> 16
$ws.Range("K7").Value = 36  # Total: 35 -> 36

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 55  # Aggravated Assault: 53 -> 55
$ws.Range("K7").Value = 221  # Total: 219 -> 221

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K5").Value = 6  # Robbery: 5 -> 6
$ws.Range("K6").Value = 10  # Total: 9 -> 10

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 13  # Aggravated Battery: 12 -> 13
$ws.Range("K7").Value = 68  # Total: 67 -> 68

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K3").Value = 4  # Aggravated Battery: 3 -> 4
$ws.Range("K7").Value = 13  # Total: 12 -> 13

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 20  # Aggravated Assault: 19 -> 20
$ws.Range("K7").Value = 58  # Total: 57 -> 58

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K6").Value = 43  # Robbery: 42 -> 43
$ws.Range("K7").Value = 92  # Total: 91 -> 92

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 22  # Aggravated Battery: 21 -> 22
$ws.Range("K7").Value = 57  # Total: 56 -> 57

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 60  # Aggravated Battery: 59 -> 60
$ws.Range("K6").Value = 36  # Robbery: 35 -> 36
$ws.Range("K7").Value = 172  # Total: 170 -> 172

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 10  # Aggravated Assault: 9 -> 10
$ws.Range("K3").Value = 14  # Aggravated Battery: 13 -> 14
$ws.Range("K7").Value = 44  # Total: 42 -> 44

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 42  # Aggravated Assault: 40 -> 42
$ws.Range("K6").Value = 52  # Robbery: 51 -> 52
$ws.Range("K7").Value = 142  # Total: 139 -> 142

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 15  # Aggravated Battery: 14 -> 15
$ws.Range("K7").Value = 47  # Total: 46 -> 47

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("K2").Value = 7  # Aggravated Assault: 6 -> 7
$ws.Range("K6").Value = 10  # Total: 9 -> 10

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 64  # Aggravated Assault: 63 -> 64
$ws.Range("K7").Value = 178  # Total: 177 -> 178

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 21  # Aggravated Assault: 20 -> 21
$ws.Range("K7").Value = 76  # Total: 75 -> 76

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 6  # Aggravated Assault: 5 -> 6
$ws.Range("K7").Value = 35  # Total: 34 -> 35

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K5").Value = 8  # Robbery: 7 -> 8
$ws.Range("K6").Value = 10  # Total: 9 -> 10

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 15  # Aggravated Assault: 14 -> 15
$ws.Range("K7").Value = 51  # Total: 50 -> 51

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K3").Value = 9  # Aggravated Battery: 8 -> 9
$ws.Range("K6").Value = 37  # Robbery: 35 -> 37
$ws.Range("K7").Value = 58  # Total: 55 -> 58

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K5").Value = 6  # Robbery: 5 -> 6
$ws.Range("K6").Value = 44  # Total: 43 -> 44

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 9  # Aggravated Assault: 8 -> 9
$ws.Range("K7").Value = 23  # Total: 22 -> 23

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K5").Value = 2  # Homicide: 1 -> 2
$ws.Range("K7").Value = 75  # Total: 74 -> 75

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 2  # Aggravated Battery: 3 -> 2
$ws.Range("K7").Value = 14  # Total: 15 -> 14

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K3").Value = 17  # Aggravated Battery: 16 -> 17
$ws.Range("K7").Value = 43  # Total: 42 -> 43

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 10  # Aggravated Assault: 9 -> 10
$ws.Range("K7").Value = 62  # Total: 61 -> 62

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 116  # Aggravated Assault: 115 -> 116
$ws.Range("K3").Value = 107  # Aggravated Battery: 105 -> 107
$ws.Range("K6").Value = 79  # Robbery: 77 -> 79
$ws.Range("K7").Value = 324  # Total: 319 -> 324

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 5  # Aggravated Battery: 4 -> 5
$ws.Range("K7").Value = 17  # Total: 16 -> 17

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K5").Value = 15  # Robbery: 13 -> 15
$ws.Range("K6").Value = 26  # Total: 24 -> 26

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K5").Value = 6  # Robbery: 5 -> 6
$ws.Range("K6").Value = 10  # Total: 9 -> 10
